# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet: new case/death/recovered counts for the
# affected countries, the refreshed "Datos actualizados" timestamp, and four
# rows whose country ended up in a different position in the source feed
# (Barein/Rumania/Suiza, Bulgaria/Finlandia, Zimbabue/Uganda/..., and
# Groenlandia/Islas Malvinas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 23:39"
$ws.Range("B4").Value = 3471694
$ws.Range("C4").Value = 57699
$ws.Range("D4").Value = 1542163
$ws.Range("E4").Value = 1791396
$ws.Range("G4").Value = 353
$ws.Range("H4").Value = 138135
$ws.Range("B5").Value = 1884967
$ws.Range("C5").Value = 18791
$ws.Range("E5").Value = 598622
$ws.Range("G5").Value = 682
$ws.Range("H5").Value = 72833
$ws.Range("B8").Value = 330123
$ws.Range("C8").Value = 3797
$ws.Range("D8").Value = 221008
$ws.Range("E8").Value = 97061
$ws.Range("G8").Value = 184
$ws.Range("H8").Value = 12054
$ws.Range("B13").Value = 287796
$ws.Range("C13").Value = 11554
$ws.Range("D13").Value = 138241
$ws.Range("E13").Value = 145383
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 4172
$ws.Range("B19").Value = 200431
$ws.Range("C19").Value = 481
$ws.Range("E19").Value = 6192
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 9139
$ws.Range("D21").Value = 78597
$ws.Range("E21").Value = 63751
$ws.Range("B27").Value = 83001
$ws.Range("C27").Value = 931
$ws.Range("D27").Value = 24975
$ws.Range("E27").Value = 54091
$ws.Range("G27").Value = 77
$ws.Range("H27").Value = 3935
$ws.Range("B31").Value = 68459
$ws.Range("C31").Value = 589
$ws.Range("D31").Value = 30369
$ws.Range("E31").Value = 33027
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 5063
$ws.Range("A49").Value = "Barein"
$ws.Range("B49").Value = 33476
$ws.Range("C49").Value = 535
$ws.Range("D49").Value = 29099
$ws.Range("E49").Value = 4268
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 109
$ws.Range("A50").Value = "Rumania"
$ws.Range("B50").Value = 32948
$ws.Range("C50").Value = 413
$ws.Range("D50").Value = 21692
$ws.Range("E50").Value = 9355
$ws.Range("G50").Value = 17
$ws.Range("H50").Value = 1901
$ws.Range("A51").Value = "Suiza"
$ws.Range("B51").Value = 32946
$ws.Range("C51").Value = 63
$ws.Range("D51").Value = 29600
$ws.Range("E51").Value = 1378
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 1968
$ws.Range("B71").Value = 12872
$ws.Range("C71").Value = 106
$ws.Range("D71").Value = 6810
$ws.Range("E71").Value = 5978
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 7411
$ws.Range("C85").Value = 159
$ws.Range("D85").Value = 3517
$ws.Range("E85").Value = 3618
$ws.Range("G85").Value = 8
$ws.Range("H85").Value = 276
$ws.Range("A86").Value = "Finlandia"
$ws.Range("B86").Value = 7295
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 6800
$ws.Range("E86").Value = 166
$ws.Range("H86").Value = 329
$ws.Range("E106").Value = 1662
$ws.Range("G106").Value = 3
$ws.Range("H106").Value = 25
$ws.Range("B132").Value = 1378
$ws.Range("C132").Value = 41
$ws.Range("D132").Value = 710
$ws.Range("E132").Value = 664
$ws.Range("A140").Value = "Zimbabue"
$ws.Range("B140").Value = 1034
$ws.Range("C140").Value = 49
$ws.Range("D140").Value = 343
$ws.Range("E140").Value = 672
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 19
$ws.Range("A141").Value = "Uganda"
$ws.Range("B141").Value = 1029
$ws.Range("C141").Value = 4
$ws.Range("D141").Value = 977
$ws.Range("E141").Value = 52
$ws.Range("H141").Value = 0
$ws.Range("A142").Value = "Liberia"
$ws.Range("B142").Value = 1024
$ws.Range("C142").Value = 14
$ws.Range("D142").Value = 439
$ws.Range("E142").Value = 534
$ws.Range("H142").Value = 51
$ws.Range("A143").Value = "Republica de Chipre"
$ws.Range("B143").Value = 1022
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 839
$ws.Range("E143").Value = 164
$ws.Range("H143").Value = 19
$ws.Range("A144").Value = "Georgia"
$ws.Range("B144").Value = 995
$ws.Range("C144").Value = 9
$ws.Range("D144").Value = 857
$ws.Range("E144").Value = 123
$ws.Range("H144").Value = 15
$ws.Range("A145").Value = "Uruguay"
$ws.Range("B145").Value = 987
$ws.Range("D145").Value = 896
$ws.Range("E145").Value = 60
$ws.Range("H145").Value = 31
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
